# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8862
    3  = 8280
    4  = 145
    8  = 148
    9  = 159
    12 = 751
    13 = 207
    14 = 5355
    15 = 3
    18 = 13
    22 = 168
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
